$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the data cells in row 2 that hold values referencing removed shared strings,
# but keep I2 / J2 (they retain their style but lose their content).
$ws.Range("A2").Value = $null
$ws.Range("B2").Value = $null
$ws.Range("C2").Value = $null
$ws.Range("E2").Value = $null
$ws.Range("I2").Value = $null
$ws.Range("J2").Value = $null
$ws.Range("M2").Value = $null
$ws.Range("N2").Value = $null
$ws.Range("O2").Value = $null

# Update the selection shown when the workbook is reopened.
$ws.Range("A2:R6").Select()
